# Update January 2017 dividend figures on the "Yearly" sheet.
$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("L3").Value = 50.57
$wsYearly.Range("M3").Value = 12.74
$wsYearly.Range("N3").Value = 5.26

# Mirror the same edits on the "All Time" sheet (row 8 = 2017).
# F8 already pulls from Yearly!L3 via formula, so only G8/H8 need updating.
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Range("G8").Value = 12.74
$wsAllTime.Range("H8").Value = 5.26

# Reflect the final selection/active sheet state recorded by Excel.
$wsYearly.Range("N4").Select()
$wsAllTime.Range("K16").Select()
$wsAllTime.Activate()
